$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new record was logged for 2026/01/20 (火) 13:00, ranking 17.
# It belongs right before the existing run of 2026/12/29 rows, so insert
# a fresh row at 662 and push everything below it down by one.
$ws.Rows.Item(662).Insert()

$newRow = $ws.Rows.Item(662)
$newRow.Cells.Item(1, 1).NumberFormat = "@"
$newRow.Cells.Item(1, 1).Value = "2026/01/20"
$newRow.Cells.Item(1, 1).Style = "Normal"
$newRow.Cells.Item(1, 2).Value = "火"
$newRow.Cells.Item(1, 3).Value = 13
$newRow.Cells.Item(1, 4).Value = 17
